$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Meta/Venda values for rows 2-16 (row 1 is the header and is unchanged)
$rows = @(
    @(5000, 6000),   # row 2
    @(5000, 5000),   # row 3
    @(5000, 5000),   # row 4
    @(5000, 5000),   # row 5
    @(1000, 1000),   # row 6
    @(5000, 5000),   # row 7
    @(5000, 5000),   # row 8
    @(9999, 9999),   # row 9  (new)
    @(5000, 5000),   # row 10 (new)
    @(5000, 5000),   # row 11 (new)
    @(5000, 5000),   # row 12 (new)
    @(5000, 5000),   # row 13 (new)
    @(5000, 5000),   # row 14 (new)
    @(5000, 5000),   # row 15 (new)
    @(5000, 5000)    # row 16 (new)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $pair = $rows[$i]
    $ws.Range("A$r").Value = $pair[0]
    $ws.Range("B$r").Value = $pair[1]
}
